$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.2089093701996928
    "C2" = 0.5407066052227343
    "J2" = 0.009216589861751152
    "O2" = 0.001536098310291859
    "P2" = 0.1397849462365591
    "S2" = 0.09984639016897082
    "B3" = 0.0108695652173913
    "C3" = 0.01902173913043478
    "J3" = 0.02445652173913044
    "P3" = 0.7309782608695652
    "S3" = 0.2146739130434783
    "J4" = 0.04545454545454546
    "O4" = 0.01136363636363636
    "P4" = 0.7045454545454546
    "S4" = 0.2386363636363636
    "P5" = 0.5
    "S5" = 0.5
    "B6" = 0.09583333333333334
    "D6" = 0.0125
    "F6" = 0.08958333333333333
    "J6" = 0.1979166666666667
    "O6" = 0.04375
    "Q6" = 0.1520833333333333
    "R6" = 0.075
    "S6" = 0.3333333333333333
    "B7" = 0.1266490765171504
    "D7" = 0.02902374670184697
    "F7" = 0.06860158311345646
    "J7" = 0.09762532981530343
    "O7" = 0.01319261213720317
    "Q7" = 0.1794195250659631
    "R7" = 0.07387862796833773
    "S7" = 0.4116094986807388
    "B8" = 0.1002132196162047
    "D8" = 0.01492537313432836
    "F8" = 0.07142857142857142
    "J8" = 0.1162046908315565
    "O8" = 0.02345415778251599
    "Q8" = 0.1428571428571428
    "R8" = 0.07889125799573561
    "S8" = 0.4520255863539446
    "B9" = 0.1253822629969419
    "D9" = 0.01529051987767584
    "E9" = 0.003058103975535168
    "F9" = 0.08868501529051988
    "J9" = 0.1039755351681957
    "O9" = 0.02140672782874618
    "Q9" = 0.1345565749235474
    "R9" = 0.08868501529051988
    "S9" = 0.418960244648318
    "B10" = 0.1189351653069987
    "D10" = 0.0240446543580936
    "E10" = 0.0008587376556462001
    "F10" = 0.07170459424645771
    "J10" = 0.1155002146844139
    "O10" = 0.0206097037355088
    "Q10" = 0.1923572348647488
    "R10" = 0.0738514383855732
    "S10" = 0.382138256762559
    "G11" = 0.1414634146341463
    "J11" = 0.09918699186991869
    "K11" = 0.2260162601626016
    "L11" = 0.5186991869918699
    "S11" = 0.01463414634146342
    "G12" = 0.7243401759530792
    "J12" = 0.1906158357771261
    "K12" = 0.01173020527859238
    "L12" = 0.04105571847507331
    "S12" = 0.03225806451612903
    "G13" = 0.6790123456790124
    "J13" = 0.2592592592592592
    "S13" = 0.06172839506172839
    "G14" = 0.8333333333333334
    "J14" = 0.1666666666666667
    "F15" = 0.02644230769230769
    "H15" = 0.1466346153846154
    "I15" = 0.0625
    "J15" = 0.3365384615384616
    "K15" = 0.06971153846153846
    "M15" = 0.01442307692307692
    "N15" = 0.002403846153846154
    "O15" = 0.06009615384615385
    "S15" = 0.28125
    "F16" = 0.01699029126213592
    "H16" = 0.2257281553398058
    "I16" = 0.08009708737864078
    "J16" = 0.3519417475728155
    "K16" = 0.1213592233009709
    "M16" = 0.01456310679611651
    "N16" = 0.002427184466019417
    "O16" = 0.03155339805825243
    "S16" = 0.1553398058252427
    "F17" = 0.01849405548216645
    "H17" = 0.2113606340819023
    "I17" = 0.07529722589167767
    "J17" = 0.4187582562747688
    "K17" = 0.1056803170409511
    "M17" = 0.01321003963011889
    "N17" = 0.001321003963011889
    "O17" = 0.05151915455746367
    "S17" = 0.1043593130779392
    "F18" = 0.0208955223880597
    "H18" = 0.2388059701492537
    "I18" = 0.0417910447761194
    "J18" = 0.417910447761194
    "K18" = 0.08358208955223881
    "M18" = 0.01791044776119403
    "N18" = 0.002985074626865672
    "O18" = 0.0626865671641791
    "S18" = 0.1134328358208955
    "F19" = 0.02106518282988871
    "H19" = 0.2209856915739269
    "I19" = 0.07869634340222575
    "J19" = 0.3608903020667726
    "K19" = 0.1104928457869634
    "M19" = 0.02384737678855326
    "N19" = 0.000794912559618442
    "O19" = 0.06756756756756757
    "S19" = 0.1156597774244833
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
